# Temps de simulation sur AMD FX-8350
# Adds a second results column ("AMD FX-8350 4 GHz (Run 1)") next to the
# existing "GAB" (SPS/PSIM) timing table on Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------
# 1. Make room: insert a new blank row 1 (everything currently on rows
#    1-8 shifts down to rows 2-9, merges shift with it).
# ---------------------------------------------------------------------
$ws.Rows("1:1").Insert()

# ---------------------------------------------------------------------
# 2. New header row (row 1): "GAB" over C1:D1, "AMD FX-8350 4 GHz (Run 1)"
#    over E1:F1.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "GAB"
$ws.Range("E1").Value = "AMD FX-8350 4 GHz (Run 1)"

$ws.Range("C1:D1").Merge()
$ws.Range("E1:F1").Merge()

$ws.Range("C1:D1").HorizontalAlignment = -4108
$ws.Range("C1:D1").Borders.Item(9).LineStyle = 1

$ws.Range("E1:F1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Row 2 (was row 1): add SPS/PSIM headers over the new E/F columns too.
# ---------------------------------------------------------------------
$ws.Range("E2").Value = "SPS"
$ws.Range("F2").Value = "PSIM"
$ws.Range("E2:F2").Font.Bold = $true
$ws.Range("E2:F2").Borders.LineStyle = 1
$ws.Range("E2:F2").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 4. New AMD FX-8350 timings for rows 3-9 (E = SPS run, F = PSIM run).
# ---------------------------------------------------------------------
$amd = @{
    3 = @("5m 54s ", "7m 23s")
    4 = @("3m 08s", "1m 28s")
    5 = @("4m 11s", "1m 41s")
    6 = @("2m 33s", "1m 19s")
    7 = @("2m 49s", "2m 5s")
    8 = @("2m 48s", "1m 41s")
    9 = @("1m 20s", "16s")
}

foreach ($r in 3..9) {
    $vals = $amd[$r]
    $ws.Range("E$r").Value = $vals[0]
    $ws.Range("F$r").Value = $vals[1]
}

$ws.Range("E3:F9").Borders.LineStyle = 1
$ws.Range("E3:F9").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Column widths for the two new columns.
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 12.02
$ws.Columns("F").ColumnWidth = 13.17

# ---------------------------------------------------------------------
# 6. Selection, matching the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("C20").Select()

Write-Output "done"
